$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Valid": add a hyperlink on C2 (mirrors the existing B2 hyperlink)
# and extend the selection to B2:C2.
# ---------------------------------------------------------------------------
$wsValid = $wb.Worksheets.Item("Valid")
$wsValid.Activate()

$wsValid.Hyperlinks.Add($wsValid.Range("C2"), "mailto:ipmcloud@200#", "", "", "ipmcloud@200#") | Out-Null
$wsValid.Range("C2").Style = "Hyperlink"
$wsValid.Range("C2").Value = "ipmcloud@200#"

$wsValid.Range("B2:C2").Select()

# ---------------------------------------------------------------------------
# Sheet "Invalid": retype the wrong-credential cells while keeping the old
# hyperlink targets/display text, add a second hyperlink on C2, and move
# the selection to C4.
# ---------------------------------------------------------------------------
$wsInvalid = $wb.Worksheets.Item("Invalid")
$wsInvalid.Activate()

$wsInvalid.Range("B2").Hyperlinks.Delete()
$wsInvalid.Hyperlinks.Add($wsInvalid.Range("B2"), "mailto:test1@vipl.com", "", "", "test1@vipl.com") | Out-Null
$wsInvalid.Range("B2").Style = "Hyperlink"
$wsInvalid.Range("B2").Value = "adas"

$wsInvalid.Hyperlinks.Add($wsInvalid.Range("C2"), "mailto:ipmcloud@200#", "", "", "ipmcloud@200#") | Out-Null
$wsInvalid.Range("C2").Style = "Hyperlink"
$wsInvalid.Range("C2").Value = "asdasd"

$wsInvalid.Range("C4").Select()
